$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.547.11'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '2.616.31'
$ws.Range('E3').Value = '  +3.41%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '109.43'
$ws.Range('E5').Value = '  +2.22%  '
$ws.Range('D6').Value = '321.12'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  -1.25%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '0.536'
$ws.Range('E9').Value = '  -1.97%  '
$ws.Range('D10').Value = '39.28'
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('D11').Value = '19.62'
$ws.Range('E11').Value = '  -2.63%  '
$ws.Range('D12').Value = '0.0806'
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').Value = '7.16'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').Value = '3.022.09'
$ws.Range('E15').Value = '  +3.28%  '
$ws.Range('D16').Value = '2.609.01'
$ws.Range('E16').Value = '  +3.37%  '
$ws.Range('D17').Value = '0.854'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = '48.579.74'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').Value = '12.74'
$ws.Range('E19').Value = '  -1.87%  '
$ws.Range('B20').Value = 'ImmutableX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.90'
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '6.63'
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('D22').Value = '0.0₃0936'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '269.60'
$ws.Range('E23').Value = '  -5.92%  '
$ws.Range('D24').Value = '69.21'
$ws.Range('E24').Value = '  -2.95%  '
$ws.Range('D25').Value = '2.52'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = '25.97'
$ws.Range('E26').Value = '  -0.33%  '
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.00'
$ws.Range('E28').Value = '  +2.86%  '
$ws.Range('E29').Value = '  -0.41%  '
$ws.Range('D30').Value = '34.82'
$ws.Range('E30').Value = '  -0.91%  '
$ws.Range('D31').Value = '0.136'
$ws.Range('E31').Value = '  -5.23%  '
$ws.Range('D32').Value = '49.08'
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('D33').Value = '5.37'
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('E34').Value = '  -1.55%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').Value = '0.0789'
$ws.Range('E36').Value = '  +1.73%  '
$ws.Range('D37').Value = '4.88'
$ws.Range('E37').Value = '  +6.15%  '
$ws.Range('D38').Value = '2.01'
$ws.Range('E38').Value = '  +1.49%  '
$ws.Range('D39').Value = '3.09'
$ws.Range('E39').Value = '  +5.57%  '
$ws.Range('D40').Value = '124.57'
$ws.Range('E40').Value = '  +3.55%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.110'
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '22.25'
$ws.Range('E42').Value = '  +1.85%  '
$ws.Range('D43').Value = '2.15'
$ws.Range('E43').Value = '  -3.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0310'
$ws.Range('E44').Value = '  +2.06%  '
$ws.Range('D45').Value = '2.058.98'
$ws.Range('E45').Value = '  +2.60%  '
$ws.Range('D46').Value = '3.18'
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('E47').Value = '  +4.24%  '
$ws.Range('E48').Value = '  +3.19%  '
$ws.Range('D49').Value = '8.88'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('D50').Value = '58.32'
$ws.Range('E50').Value = '  +3.36%  '
$ws.Range('D51').Value = '5.13'
$ws.Range('E51').Value = '  -1.41%  '
